# Apply the "state csv reading done" commit:
#  - update the stored absolute path + active tab hint on the workbook
#  - on the ALLKEYS sheet, add a new column G with fillFormField(...) helper
#    formulas for every populated data row
#  - move the UI focus (active sheet / selection) onto ALLKEYS, which also
#    clears the stale tabSelected flag that used to sit on the GOA sheet

$wb = $excel.ActiveWorkbook

# ---- workbook-level bits -------------------------------------------------
$wb.FullName = "C:\Users\Nagasudhir\Documents\WRLDC\data-portal\constkeystringconfig.xlsx"

$ws = $wb.Worksheets.Item("ALLKEYS")

# Data rows on ALLKEYS (blank separator rows -- 15, 22, 29, 36, 49, 61, 77 -- excluded)
$dataRows = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,23,24,25,26,27,28,30,31,32,33,34,35,37,38,39,40,41,42,43,44,45,46,47,48,50,51,52,53,54,55,56,57,58,59,60,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92)

foreach ($r in $dataRows) {
    $aRef = "A" + $r
    $gRef = "G" + $r
    $ws.Range($gRef).Formula = "=""fillFormField('""&" + $aRef + "&""', ""&" + $aRef + "&"");"""
}

# Column G width, to match the rest of the sheet's bestFit columns
$ws.Columns.Item(7).ColumnWidth = 44.85546875

# ---- selection / view state ---------------------------------------------
$ws.Activate()
$ws.Range("B76").Select()
$ws.Range("G78:G92").Select()
